# Apply "Trade #12 closed" update to the live trading results workbook.
$wb = $excel.ActiveWorkbook

# ---- Sheet: Summary ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1200.05
$summary.Range("B4").Value = 0.04
$summary.Range("B5").Value = 0.07000000000000001
$summary.Range("B6").Value = 12
$summary.Range("B7").Value = 6
$summary.Range("B9").Value = 50

# ---- Sheet: Strategy Status (MarketMaking row = row 4) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 100.05
$status.Range("D4").Value = 12
$status.Range("E4").Value = 0.04
$status.Range("F4").Value = 0.05
$status.Range("G4").Value = 50

# ---- New trade row appended to "All Trades" and "MarketMaking" sheets ----
$newRow = @(12, "2026-02-17", "12:27:45", "MarketMaking", "DOWN", 0.9, 0.91, "CLOSED", 1.1111, 0.01, 100.05, 0, 0, 0.6, "Normal spread capture: 19600 bps", "early_exit", 0.13)
# Columns B (Date) and C (Time) hold plain text like "2026-02-17" / "12:27:45" in
# this sheet, so force text format before assigning - otherwise Excel's COM
# value-setter auto-parses the date-shaped string into a date serial number.
$textCols = @(2, 3)

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    for ($col = 1; $col -le $newRow.Length; $col++) {
        $cell = $ws.Cells.Item(13, $col)
        if ($textCols -contains $col) {
            $cell.NumberFormat = "@"
            $cell.Value = $newRow[$col - 1]
            # Drop the temporary text format again so the cell ends up with no
            # explicit style, matching the rest of the sheet's plain cells.
            $cell.ClearFormats()
        } else {
            $cell.Value = $newRow[$col - 1]
        }
    }
}
